$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows ---
# E11: "Mockup Fini" -> "Zoning Fini"
$ws.Range("E11").Value = "Zoning Fini"

# E15: "Poursuite du CDC, charte graph, moodboard, arbo" -> same text, now at new shared-string index
$ws.Range("E15").Value = "Poursuite du CDC, charte graph, moodboard, arbo"

# --- Add new rows 17,19,21,23,25, copying the date style from C9 ---
$ws.Range("C9").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C21").PasteSpecial(-4122)
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C25").PasteSpecial(-4122)

$ws.Range("C17").Value = 45255
$ws.Range("E17").Value = "Moodboard/arborescence/wireframe terminés"

$ws.Range("C19").Value = 45256
$ws.Range("E19").Value = "Début du prototype"

$ws.Range("C21").Value = 45258
$ws.Range("E21").Value = "Poursuite proto"

$ws.Range("C23").Value = 45259
$ws.Range("E23").Value = "Prototype terminé"

$ws.Range("C25").Value = 45263
$ws.Range("E25").Value = "Diagramme de gant / mockup"

# --- Selection / view state ---
$ws.Range("H33").Select()
